$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 57.52065133333334
$ws.Range("H2").Value = 172.561954
$ws.Range("I2").Value = 0.1828443315107865
$ws.Range("J2").Value = 0.1855832454108249
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1613523333333333
$ws.Range("N2").Value = 0.484057
$ws.Range("O2").Value = 0.2926766298022186
$ws.Range("P2").Value = 0.3782017374917083
$ws.Range("Q2").Value = 9.281091307486447
$ws.Range("R2").Value = 83.52982176737801
$ws.Range("S2").Value = 0.0535142627250166
$ws.Range("T2").Value = 0.07018790586372407

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 57.52065133333334
$ws.Range("H3").Value = 172.561954
$ws.Range("I3").Value = 0.1828443315107865
$ws.Range("J3").Value = 0.1855832454108249
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01594066666666667
$ws.Range("N3").Value = 0.047822
$ws.Range("O3").Value = 0.02891473894686308
$ws.Range("P3").Value = 0.03736411928828315
$ws.Range("Q3").Value = 0.9169175293542224
$ws.Range("R3").Value = 8.252257764188002
$ws.Range("S3").Value = 0.005286896113548082
$ws.Range("T3").Value = 0.006934154519436787

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 57.52065133333334
$ws.Range("H4").Value = 172.561954
$ws.Range("I4").Value = 0.1828443315107865
$ws.Range("J4").Value = 0.1855832454108249
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3740059999999999
$ws.Range("N4").Value = 0.7480119999999999
$ws.Range("O4").Value = 0.6784086312509182
$ws.Range("P4").Value = 0.5844341432200085
$ws.Range("Q4").Value = 21.51306872257467
$ws.Range("R4").Value = 129.078412335448
$ws.Range("S4").Value = 0.1240431726722218
$ws.Range("T4").Value = 0.108461185027664

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 243.0020346666667
$ws.Range("H5").Value = 729.006104
$ws.Range("I5").Value = 0.7724450880589986
$ws.Range("J5").Value = 0.7840159175794992
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1613523333333333
$ws.Range("N5").Value = 0.484057
$ws.Range("O5").Value = 0.2926766298022186
$ws.Range("P5").Value = 0.3782017374917083
$ws.Range("Q5").Value = 39.20894529821422
$ws.Range("R5").Value = 352.880507683928
$ws.Range("S5").Value = 0.2260766250803857
$ws.Range("T5").Value = 0.2965161822497226

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 243.0020346666667
$ws.Range("H6").Value = 729.006104
$ws.Range("I6").Value = 0.7724450880589986
$ws.Range("J6").Value = 0.7840159175794992
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01594066666666667
$ws.Range("N6").Value = 0.047822
$ws.Range("O6").Value = 0.02891473894686308
$ws.Range("P6").Value = 0.03736411928828315
$ws.Range("Q6").Value = 3.873614433943112
$ws.Range("R6").Value = 34.862529905488
$ws.Range("S6").Value = 0.02233504807201261
$ws.Range("T6").Value = 0.02929406426835318

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 243.0020346666667
$ws.Range("H7").Value = 729.006104
$ws.Range("I7").Value = 0.7724450880589986
$ws.Range("J7").Value = 0.7840159175794992
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3740059999999999
$ws.Range("N7").Value = 0.7480119999999999
$ws.Range("O7").Value = 0.6784086312509182
$ws.Range("P7").Value = 0.5844341432200085
$ws.Range("Q7").Value = 90.88421897754132
$ws.Range("R7").Value = 545.305313865248
$ws.Range("S7").Value = 0.5240334149066003
$ws.Range("T7").Value = 0.4582056710614234

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1029616666666667
$ws.Range("H8").Value = 0.308885
$ws.Range("I8").Value = 0.000327290402255814
$ws.Range("J8").Value = 0.0003321930438891683
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1613523333333333
$ws.Range("N8").Value = 0.484057
$ws.Range("O8").Value = 0.2926766298022186
$ws.Range("P8").Value = 0.3782017374917083
$ws.Range("Q8").Value = 0.01661310516055556
$ws.Range("R8").Value = 0.149517946445
$ws.Range("S8").Value = 0.00009579025189884411
$ws.Range("T8").Value = 0.0001256359863815428

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1029616666666667
$ws.Range("H9").Value = 0.308885
$ws.Range("I9").Value = 0.000327290402255814
$ws.Range("J9").Value = 0.0003321930438891683
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01594066666666667
$ws.Range("N9").Value = 0.047822
$ws.Range("O9").Value = 0.02891473894686308
$ws.Range("P9").Value = 0.03736411928828315
$ws.Range("Q9").Value = 0.001641277607777778
$ws.Range("R9").Value = 0.01477149847
$ws.Range("S9").Value = 0.000009463516541040671
$ws.Range("T9").Value = 0.00001241210051861277

$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.1029616666666667
$ws.Range("H10").Value = 0.308885
$ws.Range("I10").Value = 0.000327290402255814
$ws.Range("J10").Value = 0.0003321930438891683
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3740059999999999
$ws.Range("N10").Value = 0.7480119999999999
$ws.Range("O10").Value = 0.6784086312509182
$ws.Range("P10").Value = 0.5844341432200085
$ws.Range("Q10").Value = 0.03850828110333333
$ws.Range("R10").Value = 0.23104968662
$ws.Range("S10").Value = 0.0002220366338159292
$ws.Range("T10").Value = 0.0001941449569890128

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 13.928462
$ws.Range("H11").Value = 27.856924
$ws.Range("I11").Value = 0.04427523444762439
$ws.Range("J11").Value = 0.02995896976851976
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1613523333333333
$ws.Range("N11").Value = 0.484057
$ws.Range("O11").Value = 0.2926766298022186
$ws.Range("P11").Value = 0.3782017374917083
$ws.Range("Q11").Value = 2.247389843444667
$ws.Range("R11").Value = 13.484339060668
$ws.Range("S11").Value = 0.0129583264018338
$ws.Range("T11").Value = 0.01133053441991574

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 13.928462
$ws.Range("H12").Value = 27.856924
$ws.Range("I12").Value = 0.04427523444762439
$ws.Range("J12").Value = 0.02995896976851976
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.01594066666666667
$ws.Range("N12").Value = 0.047822
$ws.Range("O12").Value = 0.02891473894686308
$ws.Range("P12").Value = 0.03736411928828315
$ws.Range("Q12").Value = 0.2220289699213334
$ws.Range("R12").Value = 1.332173819528
$ws.Range("S12").Value = 0.001280206845864219
$ws.Range("T12").Value = 0.001119390520185041

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 13.928462
$ws.Range("H13").Value = 27.856924
$ws.Range("I13").Value = 0.04427523444762439
$ws.Range("J13").Value = 0.02995896976851976
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.3740059999999999
$ws.Range("N13").Value = 0.7480119999999999
$ws.Range("O13").Value = 0.6784086312509182
$ws.Range("P13").Value = 0.5844341432200085
$ws.Range("Q13").Value = 5.209328358771999
$ws.Range("R13").Value = 20.837313435088
$ws.Range("S13").Value = 0.03003670119992636
$ws.Range("T13").Value = 0.01750904482841898

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.033993
$ws.Range("H14").Value = 0.101979
$ws.Range("I14").Value = 0.0001080555803345765
$ws.Range("J14").Value = 0.0001096741972668582
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1613523333333333
$ws.Range("N14").Value = 0.484057
$ws.Range("O14").Value = 0.2926766298022186
$ws.Range("P14").Value = 0.3782017374917083
$ws.Range("Q14").Value = 0.005484849867000001
$ws.Range("R14").Value = 0.049363648803
$ws.Range("S14").Value = 0.00003162534308364674
$ws.Range("T14").Value = 0.00004147897196433414

$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.033993
$ws.Range("H15").Value = 0.101979
$ws.Range("I15").Value = 0.0001080555803345765
$ws.Range("J15").Value = 0.0001096741972668582
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.01594066666666667
$ws.Range("N15").Value = 0.047822
$ws.Range("O15").Value = 0.02891473894686308
$ws.Range("P15").Value = 0.03736411928828315
$ws.Range("Q15").Value = 0.0005418710820000001
$ws.Range("R15").Value = 0.004876839738000001
$ws.Range("S15").Value = 0.000003124398897126071
$ws.Range("T15").Value = 0.000004097879789525587

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.033993
$ws.Range("H16").Value = 0.101979
$ws.Range("I16").Value = 0.0001080555803345765
$ws.Range("J16").Value = 0.0001096741972668582
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.3740059999999999
$ws.Range("N16").Value = 0.7480119999999999
$ws.Range("O16").Value = 0.6784086312509182
$ws.Range("P16").Value = 0.5844341432200085
$ws.Range("Q16").Value = 0.012713585958
$ws.Range("R16").Value = 0.07628151574799999
$ws.Range("S16").Value = 0.00007330583835380367
$ws.Range("T16").Value = 0.00006409734551299847

Write-Output "applied"
